$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# Set column D ("A") values for the rows that were previously blank
$ws.Range("D2").Value = "A"
$ws.Range("D5").Value = "A"
$ws.Range("D6").Value = "A"
$ws.Range("D7").Value = "A"
$ws.Range("D8").Value = "A"
$ws.Range("D11").Value = "A"

# Clear the "Pass" markers from column J (removing the "Pass"/"Fail" lookup entries)
$ws.Range("J2").ClearContents()
$ws.Range("J5").ClearContents()
$ws.Range("J6").ClearContents()
$ws.Range("J7").ClearContents()
$ws.Range("J8").ClearContents()
$ws.Range("J9").ClearContents()
$ws.Range("J10").ClearContents()
$ws.Range("J11").ClearContents()

# Update the selected cell / scroll position on the TestCases sheet
$ws.Range("D2").Select()
